# "added function to run custom query on database"
# - update a couple of existing label cells (I7/J7, I20/J20)
# - append 6 new data rows (23-28) following the existing pattern

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows -------------------------------------------------
$ws.Range("I7").Value  = " Storage"
$ws.Range("J7").Value  = " N/A"

$ws.Range("I20").Value = " Storage"
$ws.Range("J20").Value = " N/A"

# --- Append new rows 23-28 -------------------------------------------------
# Every new cell holds text that looks numeric ("123".."128", " 24".." 29"),
# so force a text number format on that block first - otherwise Excel would
# store the values as real numbers instead of text, same as it would for a
# live user typing digits into a cell that isn't pre-formatted as Text.
$ws.Range("A23:J28").NumberFormat = "@"

for ($i = 23; $i -le 28; $i++) {
    $idVal = 100 + $i
    $numVal = $i + 1

    $ws.Cells.Item($i, 1).Value2 = $idVal

    for ($col = 2; $col -le 10; $col++) {
        $ws.Cells.Item($i, $col).Value2 = " " + $numVal
    }
}
